$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels (columns G, H, I) ---
# Insertion order controls the shared-string table order, so we add
# G1's string first, then I1's, then H1's, to match the expected
# shared string indices (12=we_diameter_mm, 13=cell_capacity_mah_cm2,
# 14=we_area_cm2).
$ws.Range("G1").Value = "we_diameter_mm"
$ws.Range("I1").Value = "cell_capacity_mah_cm2"
$ws.Range("H1").Value = "we_area_cm2"

# --- New data: we_diameter_mm column (G) is a constant 14 mm for every
# data row (skipping the "blank" separator rows 2, 6 and 10) ---
$ws.Range("G3").Value = 14
$ws.Range("G4:G13").Value = 14

# --- New formula columns: we_area_cm2 (H) and cell_capacity_mah_cm2 (I) ---
# Row 3 is entered as standalone (non-shared) formulas, matching the
# first row of each new column being typed directly.
$ws.Range("H3").Formula = "=PI() * ( G3 / 2 / 10 )^2"
$ws.Range("I3").Formula = "=F3/H3"

# Rows 4-13 are filled as one contiguous block (creating a shared
# formula group), then the separator rows are cleared below.
$ws.Range("H4:H13").Formula = "=PI() * ( G4 / 2 / 10 )^2"
$ws.Range("I4:I13").Formula = "=F4/H4"

# Rows 6 and 10 are "blank" separator rows - they should not carry the
# new we_diameter_mm / we_area_cm2 / cell_capacity_mah_cm2 values.
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()

# --- Final selection state ---
$ws.Range("D13").Select()
